$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark (it currently sits in the
#        title paragraph, right after "CMP73010"). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Locate the paragraph that ends with
#        ">>>  your stuff after this line >>>" so we can add a new
#        paragraph right after it. ---
$anchor = $d.Content
$anchor.Find.Execute("stuff after this line", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorParagraph = $anchor.Paragraphs(1)

# --- 3. Insert a brand new paragraph after the anchor paragraph and give
#        it the text "Changes Made by Isuru Gamage". ---
$anchorParagraph.Range.InsertParagraphAfter()
$newParagraph = $anchorParagraph.Next()
$newParagraph.Range.Text = "Changes Made by Isuru Gamage"

# --- 4. Re-add the "_GoBack" bookmark directly after the new text, still
#        inside the new paragraph (immediately before its paragraph mark).
#        A collapsed range placed exactly at "end-of-text" sometimes lands
#        oddly right at the paragraph-mark boundary, so we instead append a
#        temporary marker character, wrap the bookmark around it, and then
#        delete the marker again -- the bookmark collapses to the correct
#        spot and survives the deletion. ---
$paraRange = $newParagraph.Range
$markerStart = $paraRange.End - 1
$paraRange.InsertAfter("X")
$markerRange = $d.Range($markerStart, $markerStart + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Range($markerStart, $markerStart + 1).Delete()
